$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - copy formatting (bold, centered, bordered) from existing header cell
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for columns I (I0) and J (IF)
$data = @(
    @(6, 9),
    @(8, 9),
    @(7, 8),
    @(4, 8),
    @(2, 5),
    @(6, 7),
    @(9, 9),
    @(9, 9),
    @(7, 9),
    @(11, 11),
    @(8, 8),
    @(7, 7),
    @(7, 9),
    @(6, 7),
    @(6, 7),
    @(6, 7),
    @(6, 8),
    @(6, 9),
    @(1, 3),
    @(5, 7),
    @(4, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
